$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z1").Value = "'27.622.90"
$ws.Range("AA1").Value = "'  -1.43%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D2:E2").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.753.96"
$ws.Range("AA1").Value = "'  -0.73%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D3:E3").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'  +0.05%  "
$ws.Range("Z1").Copy()
$ws.Range("E4").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'324.05"
$ws.Range("AA1").Value = "'  +1.10%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D5:E5").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.000"
$ws.Range("AA1").Value = "'  +0.00%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D6:E6").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.4598"
$ws.Range("AA1").Value = "'  +7.81%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D7:E7").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.3595"
$ws.Range("AA1").Value = "'  -0.62%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D8:E8").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.07531"
$ws.Range("AA1").Value = "'  +1.07%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D9:E9").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'42.21"
$ws.Range("AA1").Value = "'  -3.34%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D10:E10").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.099"
$ws.Range("AA1").Value = "'  -0.09%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D11:E11").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.000"
$ws.Range("AA1").Value = "'  +0.08%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D12:E12").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'20.82"
$ws.Range("AA1").Value = "'  -1.51%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D13:E13").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'6.023"
$ws.Range("AA1").Value = "'  -0.84%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D14:E14").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'7.113"
$ws.Range("AA1").Value = "'  -2.94%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D15:E15").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.753.89"
$ws.Range("AA1").Value = "'  -1.97%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D16:E16").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'92.70"
$ws.Range("AA1").Value = "'  +1.61%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D17:E17").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.00001067"
$ws.Range("AA1").Value = "'  +0.82%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D18:E18").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.06421"
$ws.Range("AA1").Value = "'  +0.51%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D19:E19").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.001"
$ws.Range("AA1").Value = "'  +0.07%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D20:E20").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'16.78"
$ws.Range("AA1").Value = "'  -1.58%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D21:E21").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'5.831"
$ws.Range("AA1").Value = "'  -2.59%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D22:E22").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'27.667.24"
$ws.Range("AA1").Value = "'  -1.27%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D23:E23").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'  -0.14%  "
$ws.Range("Z1").Copy()
$ws.Range("E24").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'2.113"
$ws.Range("AA1").Value = "'  -0.89%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D25:E25").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'164.15"
$ws.Range("AA1").Value = "'  +4.19%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D26:E26").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'20.46"
$ws.Range("AA1").Value = "'  +1.22%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D27:E27").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.954.89"
$ws.Range("AA1").Value = "'  -1.63%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D28:E28").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'2.093"
$ws.Range("AA1").Value = "'  -2.09%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D29:E29").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'127.06"
$ws.Range("AA1").Value = "'  +1.73%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D30:E30").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.076"
$ws.Range("AA1").Value = "'  -6.99%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D31:E31").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.09218"
$ws.Range("AA1").Value = "'  +3.84%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D32:E32").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'3.670"
$ws.Range("AA1").Value = "'  +1.99%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D33:E33").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'5.550"
$ws.Range("AA1").Value = "'  -1.53%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D34:E34").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'11.96"
$ws.Range("AA1").Value = "'  -4.72%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D35:E35").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.02301"
$ws.Range("AA1").Value = "'  -0.63%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D36:E36").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.2101"
$ws.Range("AA1").Value = "'  -0.16%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D37:E37").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.06046"
$ws.Range("AA1").Value = "'  +0.16%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D38:E38").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.6384"
$ws.Range("AA1").Value = "'  +0.19%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D39:E39").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'4.983"
$ws.Range("AA1").Value = "'  -0.88%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D40:E40").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.204"
$ws.Range("AA1").Value = "'  +1.31%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D41:E41").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.384"
$ws.Range("AA1").Value = "'  -0.79%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D42:E42").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'7.822"
$ws.Range("AA1").Value = "'  -0.32%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D43:E43").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'13.23"
$ws.Range("AA1").Value = "'  -1.07%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D44:E44").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.5917"
$ws.Range("AA1").Value = "'  -0.08%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D45:E45").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'3.711"
$ws.Range("AA1").Value = "'  +0.67%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D46:E46").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'123.21"
$ws.Range("AA1").Value = "'  +0.60%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D47:E47").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.961"
$ws.Range("AA1").Value = "'  -2.37%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D48:E48").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'1.148"
$ws.Range("AA1").Value = "'  -3.24%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D49:E49").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'0.06859"
$ws.Range("AA1").Value = "'  -0.11%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D50:E50").PasteSpecial(-4163)
$ws.Range("Z1").Value = "'72.27"
$ws.Range("AA1").Value = "'  -2.35%  "
$ws.Range("Z1:AA1").Copy()
$ws.Range("D51:E51").PasteSpecial(-4163)

$ws.Range("Z1:AA1").Clear()
$excel.CutCopyMode = 0
